$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.881.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "'3.782.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.85%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'421.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").Value = "'133.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").Value = "'3.766.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.64%  "

$ws.Range("D8").Value = "'0.652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").Value = "'0.189"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.97%  "

$ws.Range("D12").Value = "'0.0000431"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +63.08%  "

$ws.Range("D13").Value = "'43.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Value = "'10.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.84%  "

$ws.Range("D15").Value = "'4.366.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.34%  "

$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "'3.771.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.00%  "

$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D19").Value = "'13.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("E20").Value = "  +3.62%  "

$ws.Range("D21").Value = "'67.904.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.75%  "

$ws.Range("D22").Value = "'451.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'15.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.98%  "

$ws.Range("D24").Value = "'90.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "

$ws.Range("E25").Value = "  -4.18%  "

$ws.Range("D26").Value = "'38.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.69%  "

$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("E28").Value = "  +2.32%  "

$ws.Range("D29").Value = "'5.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.34%  "

$ws.Range("E30").Value = "  +5.96%  "

$ws.Range("D31").Value = "'12.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.55%  "

$ws.Range("D32").Value = "'2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.35%  "

$ws.Range("D33").Value = "'7.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.24%  "

$ws.Range("E34").Value = "  +2.20%  "

$ws.Range("D35").Value = "'42.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.37%  "

$ws.Range("D36").Value = "'58.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").Value = "'0.0494"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("D39").Value = "'0.0₃0751"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("E40").Value = "  +30.54%  "

$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").Value = "'0.997"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D43").Value = "'27.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +28.93%  "

$ws.Range("D44").Value = "'3.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.63%  "

$ws.Range("D45").Value = "'2.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.26%  "

$ws.Range("D46").Value = "'148.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.65%  "

$ws.Range("D47").Value = "'3.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +23.21%  "

$ws.Range("D48").Value = "'2.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.30%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.21%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'4.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.12%  "

$ws.Range("D51").Value = "'0.309"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
